$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "42.845.63"
Set-TextValue "E2" "  +4.48%  "

Set-TextValue "D3" "2.264.63"
Set-TextValue "E3" "  +4.21%  "

Set-TextValue "E4" "  +0.15%  "

Set-TextValue "D5" "249.39"
Set-TextValue "E5" "  +0.95%  "

Set-TextValue "D6" "0.635"
Set-TextValue "E6" "  +3.37%  "

Set-TextValue "D7" "71.96"
Set-TextValue "E7" "  +8.63%  "

Set-TextValue "E8" "  +0.02%  "

Set-TextValue "D9" "0.663"
Set-TextValue "E9" "  +16.96%  "

Set-TextValue "D10" "38.89"
Set-TextValue "E10" "  +8.96%  "

Set-TextValue "D11" "0.0975"
Set-TextValue "E11" "  +5.06%  "

Set-TextValue "D12" "59.69"
Set-TextValue "E12" "  -1.99%  "

Set-TextValue "D13" "7.43"
Set-TextValue "E13" "  +8.26%  "

Set-TextValue "E14" "  +1.67%  "

Set-TextValue "D15" "2.602.27"
Set-TextValue "E15" "  +4.40%  "

Set-TextValue "D16" "14.93"
Set-TextValue "E16" "  +4.58%  "

Set-TextValue "D17" "0.882"
Set-TextValue "E17" "  +3.74%  "

Set-TextValue "D18" "2.260.05"
Set-TextValue "E18" "  +4.16%  "

Set-TextValue "D19" "42.802.58"
Set-TextValue "E19" "  +4.50%  "

Set-TextValue "D20" "0.0000100"
Set-TextValue "E20" "  +7.10%  "

Set-TextValue "D21" "6.33"
Set-TextValue "E21" "  +3.96%  "

Set-TextValue "D22" "73.16"
Set-TextValue "E22" "  +2.51%  "

Set-TextValue "D23" "235.90"
Set-TextValue "E23" "  +2.74%  "

Set-TextValue "E24" "  +1.65%  "

Set-TextValue "D25" "3.94"
Set-TextValue "E25" "  +6.86%  "

Set-TextValue "B26" "Dai"
Set-TextValue "C26" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  -0.01%  "

Set-TextValue "B27" "Cosmos"
Set-TextValue "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "11.45"
Set-TextValue "E27" "  +1.40%  "

Set-TextValue "E28" "  +0.19%  "

Set-TextValue "E29" "  -1.55%  "

Set-TextValue "B30" "Toncoin"
Set-TextValue "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D30" "2.12"
Set-TextValue "E30" "  +5.80%  "

Set-TextValue "B31" "Monero"
Set-TextValue "C31" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D31" "167.53"
Set-TextValue "E31" "  -0.89%  "

Set-TextValue "B32" "EthereumClassic"
Set-TextValue "C32" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D32" "20.99"
Set-TextValue "E32" "  +3.95%  "

Set-TextValue "D33" "6.56"
Set-TextValue "E33" "  +17.01%  "

Set-TextValue "D34" "0.127"
Set-TextValue "E34" "  +5.46%  "

Set-TextValue "D35" "0.0806"
Set-TextValue "E35" "  +6.79%  "

Set-TextValue "D36" "31.47"
Set-TextValue "E36" "  +29.31%  "

Set-TextValue "D37" "0.124"
Set-TextValue "E37" "  +2.89%  "

Set-TextValue "E38" "  +11.51%  "

Set-TextValue "D39" "4.76"
Set-TextValue "E39" "  +4.62%  "

Set-TextValue "D40" "0.0319"
Set-TextValue "E40" "  +4.12%  "

Set-TextValue "E41" "  +6.25%  "

Set-TextValue "D42" "12.78"
Set-TextValue "E42" "  +12.50%  "

Set-TextValue "D43" "5.82"
Set-TextValue "E43" "  +6.40%  "

Set-TextValue "D44" "9.24"
Set-TextValue "E44" "  +9.75%  "

Set-TextValue "D45" "62.29"
Set-TextValue "E45" "  +3.58%  "

Set-TextValue "D46" "0.203"
Set-TextValue "E46" "  +6.03%  "

Set-TextValue "D47" "4.85"
Set-TextValue "E47" "  -0.44%  "

Set-TextValue "E48" "  +3.22%  "

Set-TextValue "E49" "  +0.14%  "

Set-TextValue "D50" "1.18"
Set-TextValue "E50" "  +2.35%  "

Set-TextValue "D51" "1.20"
Set-TextValue "E51" "  +4.27%  "
